# 自动更新Excel文件 - daily "remaining days" countdown refresh.
#
# Column layout (header row 1):
#   A 行号   B 店铺名称   C 地址   D 总天(total days)
#   E 剩余(days remaining)   F 开始时间(start date, stored as plain YYYYMMDD number)
#   G/H/I 备注1/2/3
#
# Business rule applied once per day to every data row:
#   - if the start date in F is not a well-formed 8-digit YYYYMMDD number,
#     the row is left untouched (bad data - can't compute a new cycle date);
#   - otherwise, if E (days remaining) is down to 1, the cycle has just
#     finished: refill E back up to the row's total (D) and roll F forward
#     by that same total, starting a fresh cycle;
#   - otherwise just count E down by one for another day elapsed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $total = $dCell.Value()
    $remaining = $eCell.Value()
    $startDate = $fCell.Value()

    if ($null -eq $remaining -or $null -eq $startDate) {
        continue
    }

    # Validate F as an 8-digit YYYYMMDD integer; skip malformed dates untouched.
    $startDateText = [string][int64]$startDate
    if ($startDateText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $eCell.Value = $total
        $fCell.Value = $startDate + $total
    } else {
        $eCell.Value = $remaining - 1
    }
}
